$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2024-03-07 Thursday" "2024-03-08 Friday"

Replace-Text "124×9=" "470×9="
Replace-Text "509×7=" "676×8="
Replace-Text "160×3=" "786×5="
Replace-Text "330×8=" "290×6="
Replace-Text "639×2=" "649×8="

Replace-Text "888×9=" "156×5="
Replace-Text "426×8=" "491×6="
Replace-Text "854×6=" "795×5="
Replace-Text "882×9=" "992×5="
Replace-Text "598×9=" "290×8="

Replace-Text "263×4=" "498×5="
Replace-Text "648×9=" "670×3="
Replace-Text "224×4=" "730×7="
Replace-Text "549×2=" "467×2="
Replace-Text "370×7=" "206×8="

Replace-Text "934×6=" "546×5="
Replace-Text "541×9=" "748×9="
Replace-Text "537×8=" "972×6="
Replace-Text "673×3=" "232×5="
Replace-Text "199×7=" "869×3="

Replace-Text "716×9=" "120×3="
Replace-Text "149×4=" "300×6="
Replace-Text "383×6=" "633×7="
Replace-Text "487×7=" "920×6="
Replace-Text "514×6=" "960×7="
